{"js": "// Revert \"my 3 commit\": remove the trailing empty paragraph and the\n// \"Have a nice day!\" paragraph that followed \"Hello, world!\", restoring\n// the document to just its first three paragraphs.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"Hello, world!\") {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  // Delete every paragraph after the anchor, always targeting the\n  // paragraph immediately following the anchor and re-resolving the\n  // paragraph collection after each deletion (deleting the document's\n  // current last paragraph twice in a row without a fresh reload can be\n  // a no-op, so we reload between deletes).\n  let remaining = items.length - anchorIndex - 1;\n  while (remaining > 0) {\n    const fresh = body.paragraphs;\n    fresh.load(\"items\");\n    await context.sync();\n    fresh.items[anchorIndex + 1].delete();\n    await context.sync();\n    remaining--;\n  }\n}\n", "ps1": "# Revert \"my 3 commit\": remove the trailing empty paragraph and the\n# \"Have a nice day!\" paragraph that followed \"Hello, world!\", restoring\n# the document to just its first three paragraphs.\n$d = $word.ActiveDocument\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -eq \"Hello, world!`r\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ne -1) {\n    # Always delete the paragraph right after the anchor, re-checking the\n    # live count each time. Deleting paragraphs one at a time in ascending\n    # order (rather than collapsing a single multi-paragraph range all the\n    # way to the document's end) avoids a no-op when the deletion would\n    # otherwise consume the document's final paragraph mark in one shot.\n    while ($d.Paragraphs.Count -gt $anchorIndex) {\n        $d.Paragraphs.Item($anchorIndex + 1).Range.Delete() | Out-Null\n    }\n}\n"}
